$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.040.95'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.131.14'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '533.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.129.84'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.465'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.57%  '
$ws.Range('E10').Value = '  +2.25%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('E12').Value = '  +3.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.660.40'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.64'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '58.068.95'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.114.62'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.13'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '355.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.167'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0877'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.95%  '
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.21'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.44'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.61%  '
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('E34').Value = '  -3.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.73'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.74'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.27'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('E39').Value = '  +4.60%  '
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.472.67'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.92%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.701'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('E43').Value = '  -4.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '37.58'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.166.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.978'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.738'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.19%  '
